$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 5 and 6 swap their full record content (observation rows for the same
# locality, re-ordered). Row 10/11 and row 12/13 likewise swap.
# Only the cells whose effective value actually changes are touched below.

# --- Row 5 ---
$ws.Range("A5").Value2 = 130981911
$ws.Range("B5").Value2 = 91829
$ws.Range("E5").Value2 = 5442
$ws.Range("F5").Value2 = 'Tallticka'
$ws.Range("G5").Value2 = 'Porodaedalea pini'
$ws.Range("H5").Value2 = '(Brot.) Murrill'
$ws.Range("J5").Value2 = ''
$ws.Range("L5").Value2 = ''
$ws.Range("M5").Value2 = ''
$ws.Range("Q5").Value2 = 437697
$ws.Range("R5").Value2 = 6792416
$ws.Range("AF5").Value2 = ''
$ws.Range("AX5").Value2 = 'Eva Löfqvist'
# --- Row 6 ---
$ws.Range("A6").Value2 = 130981914
$ws.Range("B6").Value2 = 57881
$ws.Range("E6").Value2 = 100049
$ws.Range("F6").Value2 = 'Spillkråka'
$ws.Range("G6").Value2 = 'Dryocopus martius'
$ws.Range("H6").Value2 = '(Linnaeus, 1758)'
$ws.Range("J6").Value2 = ''
$ws.Range("L6").Value2 = ''
$ws.Range("M6").Value2 = 'äldre spår'
$ws.Range("Q6").Value2 = 437688
$ws.Range("R6").Value2 = 6792409
$ws.Range("AF6").Value2 = ''
$ws.Range("AX6").Value2 = 'Eva Löfqvist, Alfhild Sehlin'
# --- Row 10 ---
$ws.Range("A10").Value2 = 130981930
$ws.Range("Q10").Value2 = 437745
$ws.Range("R10").Value2 = 6792623
# --- Row 11 ---
$ws.Range("A11").Value2 = 130981933
$ws.Range("Q11").Value2 = 437877
$ws.Range("R11").Value2 = 6792522
# --- Row 12 ---
$ws.Range("A12").Value2 = 130981909
$ws.Range("B12").Value2 = 57073
$ws.Range("D12").Value2 = 'LC'
$ws.Range("E12").Value2 = 100138
$ws.Range("F12").Value2 = 'Tjäder'
$ws.Range("G12").Value2 = 'Tetrao urogallus'
$ws.Range("H12").Value2 = 'Linnaeus, 1758'
$ws.Range("K12").Value2 = ''
$ws.Range("L12").Value2 = ''
$ws.Range("M12").Value2 = 'färsk spillning'
$ws.Range("N12").Value2 = ''
$ws.Range("Q12").Value2 = 437657
$ws.Range("R12").Value2 = 6792398
$ws.Range("AX12").Value2 = 'Eva Löfqvist, Alfhild Sehlin'
# --- Row 13 ---
$ws.Range("A13").Value2 = 130981935
$ws.Range("B13").Value2 = 79243
$ws.Range("D13").Value2 = 'NT'
$ws.Range("E13").Value2 = 6425
$ws.Range("F13").Value2 = 'Garnlav'
$ws.Range("G13").Value2 = 'Alectoria sarmentosa'
$ws.Range("H13").Value2 = '(Ach.) Ach.'
$ws.Range("K13").Value2 = ''
$ws.Range("L13").Value2 = ''
$ws.Range("M13").Value2 = ''
$ws.Range("N13").Value2 = ''
$ws.Range("Q13").Value2 = 437656
$ws.Range("R13").Value2 = 6792404
$ws.Range("AX13").Value2 = 'Eva Löfqvist'
